$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price report row was inserted above the old row 7, pushing
# the former rows 7-9 down to rows 8-10 (their data is unchanged).
$ws.Rows(7).Insert()

$ws.Range("A7").Value = 1
$ws.Range("B7").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C7").Value = "Arica y Parinacota"
$ws.Range("D7").Value = 44497
$ws.Range("E7").Value = 15
$ws.Range("F7").Value = 100112017
$ws.Range("G7").Value = "Ramas de apio"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 160
$ws.Range("K7").Value = 5000
$ws.Range("L7").Value = 6000
$ws.Range("M7").Value = 5500
$ws.Range("N7").Value = "`$/atado 7 kilos"
$ws.Range("O7").Value = "Región de Arica y Parinacota"
$ws.Range("P7").Value = 5500
$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = "Hortaliza"
